$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("float transitions")

# --- Insert a new leading column for the state "idx" ---
$ws.Columns.Item(1).Insert() | Out-Null

# --- Insert a new row before the (old) last row to host the new
#     "T_FL_FLOAT" accepting state (fixes #5375: a float that already
#     saw an exponent digit now has its own state so trailing junk
#     after the exponent is treated as a new token instead of being
#     silently folded into the number). ---
$ws.Rows.Item(8).Insert() | Out-Null

# Clone the look (font/alignment/etc.) of an existing T_FL_FLOAT cell
# onto the freshly inserted row before writing its values.
$ws.Range("I3").Copy() | Out-Null
$ws.Range("B8:I8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B8:I8").Value = "T_FL_FLOAT"

# --- Populate the new "idx" column ---
$ws.Range("A1").Value = "idx"
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7

$ws.Range("A1:A9").HorizontalAlignment = -4108
$ws.Columns.Item(1).ColumnWidth = 3.43

$ws.Range("A10").Select() | Out-Null

Write-Output "ok"
